# Update the Spanish OWASP Top 10 (2021) translations in column E
# of the "Top 10 Mapping" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top 10 Mapping")

$ws.Range("E8").Value  = "A02:2021-Fallas Criptográficas"
$ws.Range("E10").Value = "A05:2021-Configuración de Seguridad Incorrecta"
$ws.Range("E11").Value = "A06:2021-Componentes Vulnerables y Desactualizados"
$ws.Range("E14").Value = "A08:2021-Fallas en la Integridad del Software y de los Datos"
$ws.Range("E15").Value = "A09:2021-Fallas en el Registro y Monitoreo de la Seguridad*"
$ws.Range("E16").Value = "A10:2021-Falsificación de Solicitudes del Lado del Servidor (SSRF)*"

# Update the selected / active cell as recorded in the saved view state
$ws.Range("C26").Select()
